# Replaced usage of "GetAwaiter().GetResult()" with "AsyncHelper.RunSync(() => asyncMethod())"
# in the AzureStorage project -> bump its version and record it in a new "11.2.0" release column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New release column header (column G, release 11.2.0)
$ws.Range("G1").Value = "11.2.0"
$ws.Range("G1").Font.Bold = $true

# Keep the same "has a version in this release" placeholder formatting as column F
# for every project row that already participated in prior releases.
$placeholderRows = @(3,4,8,9,10,11)
foreach ($r in $placeholderRows) {
    $ws.Range("G$r").Font.Bold = $false
}

# Tardigrade.Framework.AzureStorage (row 6) gets the actual new version number.
$ws.Range("G6").Value = "6.1.0"

# Excel leaves the selection on the cell that was just edited.
$ws.Range("G6").Select()
